$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "seduce emma" -> "seduce emma - connor"
$ws.Range("A15").Value = "seduce emma - connor"

# Rename "go to the beach" -> "lay on the beach" and fix its estimated time
$ws.Range("A3").Value = "lay on the beach"
$ws.Range("C3").Value = 10

# Add new activities at the bottom of the list
$ws.Range("A20").Value = "be in the ocean"
$ws.Range("B20").Value = 9
$ws.Range("C20").Value = 4

$ws.Range("A21").Value = "swim in the pool"
$ws.Range("B21").Value = 6
$ws.Range("C21").Value = 4

$ws.Range("A22").Value = "sunbathe"
$ws.Range("B22").Value = 8
$ws.Range("C22").Value = 3

# Update the selection to reflect where the user left off
$ws.Range("A23").Select()
